$wb = $excel.ActiveWorkbook

# The "Chart" sheet holds the daily cumulative Invalid/Valid counts.
$chart = $wb.Worksheets.Item("Chart")

# Find the next empty row below the existing data (row 1 is the header).
$lastRow = $chart.Cells.Item($chart.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Append the new day's data. The date column stores plain text (not real
# dates) in this workbook, so force text entry (leading apostrophe) and then
# clear the resulting formatting so the cell keeps the default/general style
# used by every other row.
$dateCell = $chart.Cells.Item($newRow, 1)
$dateCell.Value = "'2025-11-01"
$dateCell.ClearFormats()

$chart.Cells.Item($newRow, 2).Value = 0.0
$chart.Cells.Item($newRow, 3).Value = 112.0
